# chore: adapt column header formatting to respective input file names
#
# The sheet holds a diff between two AHB ("Anwendungshandbuch") format
# versions, laid out as: <old-format columns> | diff | <new-format columns>.
# The header names used to carry generic "_old" / "_new" suffixes; they are
# renamed here to carry the concrete format-version suffixes instead
# ("_FV2210" for the left/old block, "_FV2304" for the right/new block).
# The data range is also turned into a proper Excel Table and the header
# row is frozen so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row -------------------------------------------
# Columns A:J ("..._old") -> "..._FV2210"
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value()) -replace "_old$", "_FV2210"
}

# Column K is the untouched "diff" column (index 11) - left as is.

# Columns L:U ("..._new") -> "..._FV2304"
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value()) -replace "_new$", "_FV2304"
}

# --- 2. Turn the data range into an Excel Table (ListObject) -----------
$dataRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $dataRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
